$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.184.67'
$ws.Range('E2').Value = '  +3.09%  '
$ws.Range('D3').Value = '1.580.29'
$ws.Range('E3').Value = '  +1.84%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.53'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('E6').Value = '  +5.78%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.16'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +10.00%  '
$ws.Range('E9').Value = '  +2.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0594'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.81%  '
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('D12').Value = '1.804.53'
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('D13').Value = '1.581.50'
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('D14').Value = '29.197.30'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('E15').Value = '  +2.47%  '
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.34'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '236.73'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.44'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.99'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('E23').Value = '  +2.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.08'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.92'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.13'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.59%  '
$ws.Range('E27').Value = '  +4.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.36'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0469'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.07'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('E32').Value = '  +1.60%  '
$ws.Range('D33').Value = '1.423.37'
$ws.Range('E33').Value = '  +2.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.07'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('E35').Value = '  -2.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.50'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.76'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +6.30%  '
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.530'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.38%  '
$ws.Range('E41').Value = '  +2.45%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('B43').Value = 'BitcoinSV'
$ws.Range('C43').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '52.75'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +23.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.790'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0471'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.49'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.34'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').Value = '1.716.34'
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.844'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.30'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0105'
$ws.Range('E51').Value = '  +2.07%  '
